# Auto-generated edit script applying cell value changes per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 40984.8
$ws.Range("J17").Value = 40984.8
$ws.Range("L17").Value = 122954.4
$ws.Range("N17").Value = -123290.4

$ws.Range("H47").Value = 50000
$ws.Range("I47").Value = 50000
$ws.Range("K47").Value = 50000
$ws.Range("M47").Value = -49028

$ws.Range("H86").Value = 1333.625
$ws.Range("I86").Value = 1573
$ws.Range("J86").Value = 934.6667
$ws.Range("K86").Value = 1573
$ws.Range("L86").Value = 934.6667
$ws.Range("M86").Value = -450
$ws.Range("N86").Value = -3180.6667

$ws.Range("H89").Value = 1333.625
$ws.Range("I89").Value = 1573
$ws.Range("J89").Value = 934.6667
$ws.Range("K89").Value = 7865
$ws.Range("L89").Value = 4673.3335
$ws.Range("M89").Value = -2249
$ws.Range("N89").Value = -15905.3335

$ws.Range("H100").Value = 2672.818
$ws.Range("I100").Value = 1975.5
$ws.Range("J100").Value = 3071.2856
$ws.Range("K100").Value = 1975.5
$ws.Range("L100").Value = 3071.2856
$ws.Range("M100").Value = -1434.5
$ws.Range("N100").Value = -4153.2856

$ws.Range("H132").Value = 1853.5857
$ws.Range("I132").Value = 1366.434
$ws.Range("J132").Value = 3372.353
$ws.Range("K132").Value = 4099.302
$ws.Range("L132").Value = 10117.059
$ws.Range("M132").Value = -1569.302
$ws.Range("N132").Value = -15177.059

$ws.Range("H135").Value = 758.41174
$ws.Range("I135").Value = 629.0682
$ws.Range("J135").Value = 1571.4286
$ws.Range("K135").Value = 5661.6138
$ws.Range("L135").Value = 14142.8574
$ws.Range("M135").Value = -3126.6138
$ws.Range("N135").Value = -19212.8574

$ws.Range("H137").Value = 1524.7
$ws.Range("I137").Value = 1384.9375
$ws.Range("J137").Value = 1684.4286
$ws.Range("K137").Value = 4154.8125
$ws.Range("L137").Value = 5053.2858
$ws.Range("M137").Value = -1604.8125
$ws.Range("N137").Value = -10153.2858

$ws.Range("H138").Value = 5161.3687
$ws.Range("I138").Value = 1716.2858
$ws.Range("J138").Value = 9417.058999999999
$ws.Range("K138").Value = 5148.857400000001
$ws.Range("L138").Value = 28251.177
$ws.Range("M138").Value = -8.85740000000078
$ws.Range("N138").Value = -38531.177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10119.107
$ws.Range("I32").Value = 12186.98
$ws.Range("J32").Value = 2586.1428
$ws.Range("K32").Value = 12186.98
$ws.Range("L32").Value = 2586.1428
$ws.Range("M32").Value = -11899.98
$ws.Range("N32").Value = -3160.1428

$ws.Range("H61").Value = 1149.7441
$ws.Range("I61").Value = 1073.1072
$ws.Range("J61").Value = 1292.8
$ws.Range("K61").Value = 1073.1072
$ws.Range("L61").Value = 1292.8
$ws.Range("M61").Value = -861.1071999999999
$ws.Range("N61").Value = -1716.8

$ws.Range("H74").Value = 918.67645
$ws.Range("I74").Value = 876.5862
$ws.Range("J74").Value = 1162.8
$ws.Range("K74").Value = 876.5862
$ws.Range("L74").Value = 1162.8
$ws.Range("M74").Value = -2.586199999999963
$ws.Range("N74").Value = -2910.8

$ws.Range("H77").Value = 918.67645
$ws.Range("I77").Value = 876.5862
$ws.Range("J77").Value = 1162.8
$ws.Range("K77").Value = 4382.931
$ws.Range("L77").Value = 5814
$ws.Range("M77").Value = -14.93099999999959
$ws.Range("N77").Value = -14550

$ws.Range("H102").Value = 145104.14
$ws.Range("I102").Value = 2679.75
$ws.Range("K102").Value = 2679.75
$ws.Range("M102").Value = -1057.75

$ws.Range("H122").Value = 1218.7916
$ws.Range("I122").Value = 1056.8636
$ws.Range("K122").Value = 3170.5908
$ws.Range("M122").Value = -720.5907999999999

$ws.Range("H132").Value = 1955.1041
$ws.Range("I132").Value = 1198.125
$ws.Range("J132").Value = 2712.0833
$ws.Range("K132").Value = 3594.375
$ws.Range("L132").Value = 8136.249899999999
$ws.Range("M132").Value = -1064.375
$ws.Range("N132").Value = -13196.2499

$ws.Range("H136").Value = 1149.7441
$ws.Range("I136").Value = 1073.1072
$ws.Range("J136").Value = 1292.8
$ws.Range("K136").Value = 3219.3216
$ws.Range("L136").Value = 3878.4
$ws.Range("M136").Value = -669.3215999999998
$ws.Range("N136").Value = -8978.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 22519.396
$ws.Range("I20").Value = 31396.334
$ws.Range("J20").Value = 2990.1333
$ws.Range("K20").Value = 31396.334
$ws.Range("L20").Value = 2990.1333
$ws.Range("M20").Value = -31149.334
$ws.Range("N20").Value = -3484.1333

$ws.Range("H99").Value = 1820
$ws.Range("I99").Value = 1275
$ws.Range("K99").Value = 1275
$ws.Range("M99").Value = 223

$ws.Range("H134").Value = 1904.326
$ws.Range("I134").Value = 1636.1351
$ws.Range("J134").Value = 3006.889
$ws.Range("K134").Value = 4908.4053
$ws.Range("L134").Value = 9020.667000000001
$ws.Range("M134").Value = -2373.4053
$ws.Range("N134").Value = -14090.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2394.1936
$ws.Range("I31").Value = 1568.7273
$ws.Range("J31").Value = 4412
$ws.Range("K31").Value = 1568.7273
$ws.Range("L31").Value = 4412
$ws.Range("M31").Value = -1273.7273
$ws.Range("N31").Value = -5002

$ws.Range("H34").Value = 2394.1936
$ws.Range("I34").Value = 1568.7273
$ws.Range("J34").Value = 4412
$ws.Range("K34").Value = 1568.7273
$ws.Range("L34").Value = 4412
$ws.Range("M34").Value = -1366.7273
$ws.Range("N34").Value = -4816

$ws.Range("H58").Value = 700099.2
$ws.Range("I58").Value = 1029907.44
$ws.Range("K58").Value = 1029907.44
$ws.Range("M58").Value = -1029704.44

$ws.Range("H132").Value = 387678.97
$ws.Range("I132").Value = 588837.6
$ws.Range("J132").Value = 2124.9167
$ws.Range("K132").Value = 1766512.8
$ws.Range("L132").Value = 6374.750100000001
$ws.Range("M132").Value = -1763982.8
$ws.Range("N132").Value = -11434.7501

$ws.Range("H134").Value = 1454.2766
$ws.Range("I134").Value = 1238.919
$ws.Range("J134").Value = 2251.1
$ws.Range("K134").Value = 3716.757000000001
$ws.Range("L134").Value = 6753.299999999999
$ws.Range("M134").Value = -1181.757000000001
$ws.Range("N134").Value = -11823.3

$ws.Range("H136").Value = 700099.2
$ws.Range("I136").Value = 1029907.44
$ws.Range("K136").Value = 3089722.32
$ws.Range("M136").Value = -3087172.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 42857.2
$ws.Range("J51").Value = 42857.2
$ws.Range("L51").Value = 42857.2
$ws.Range("N51").Value = -43875.2

$ws.Range("H69").Value = 205830
$ws.Range("J69").Value = 205830
$ws.Range("L69").Value = 205830
$ws.Range("N69").Value = -207328

$ws.Range("H72").Value = 205830
$ws.Range("J72").Value = 205830
$ws.Range("L72").Value = 617490
$ws.Range("N72").Value = -624978

$ws.Range("H122").Value = 3823.375
$ws.Range("I122").Value = 2877.75
$ws.Range("K122").Value = 8633.25
$ws.Range("M122").Value = -6183.25

$ws.Range("H123").Value = 11858.308
$ws.Range("J123").Value = 11858.308
$ws.Range("L123").Value = 11858.308
$ws.Range("N123").Value = -16758.308

$ws.Range("H132").Value = 1420.069
$ws.Range("I132").Value = 813.5238000000001
$ws.Range("J132").Value = 3012.25
$ws.Range("K132").Value = 2440.5714
$ws.Range("L132").Value = 9036.75
$ws.Range("M132").Value = 89.42859999999973
$ws.Range("N132").Value = -14096.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9705

$ws.Range("H46").Value = 1574.8334
$ws.Range("I46").Value = 1399.6666
$ws.Range("J46").Value = 1750
$ws.Range("K46").Value = 1399.6666
$ws.Range("L46").Value = 1750
$ws.Range("M46").Value = -1211.6666
$ws.Range("N46").Value = -2126

$ws.Range("H100").Value = 15000
$ws.Range("I100").Value = 26000
$ws.Range("K100").Value = 26000
$ws.Range("M100").Value = -25459

$ws.Range("H132").Value = 2800.7705
$ws.Range("I132").Value = 2259.42
$ws.Range("K132").Value = 6778.26
$ws.Range("M132").Value = -4248.26

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 988.5
$ws.Range("I132").Value = 660.3125
$ws.Range("J132").Value = 2301.25
$ws.Range("K132").Value = 1980.9375
$ws.Range("L132").Value = 6903.75
$ws.Range("M132").Value = 549.0625
$ws.Range("N132").Value = -11963.75

$ws.Range("H136").Value = 2152.2942
$ws.Range("I136").Value = 1747.3334
$ws.Range("K136").Value = 5242.0002
$ws.Range("M136").Value = -2692.0002
